$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 12 data rows (old rows 2-13). This shifts the remaining
# data rows (old 14-21) up to become the new rows 2-9.
$ws.Range("A2:C13").EntireRow.Delete()

# New data rows to append after the retained data (new rows 10-31).
$newRows = @(
    @(-3.904241349146918,8.71459445529258,1.42584520486685),
    @(-1.888838914724495,4.313753348130421,0.935868626374464),
    @(4.815394823367857,-4.399581175584057,-2.032744779036629),
    @(5.744059562683067,-2.208836053426431,-2.385511210331534),
    @(-0.9773453107246968,-5.042484849920636,-1.550856040532772),
    @(-5.169222395236689,-6.381682007129349,-6.139806142220171),
    @(-4.023792711588003,7.085793306277576,-4.32150216927884),
    @(1.662728214263952,4.684056318723108,-1.641669896932774),
    @(2.464888723567128,-12.57504773139953,-3.376289129257202),
    @(-0.0004609318306818891,-8.720019648625392,-0.9112097813533104),
    @(-6.561141893869462,1.067614199106516,5.132379430073989),
    @(-3.78096956106336,-6.149249487656785,3.56855486356297),
    @(4.143512133451609,-9.295757961273193,1.042706482227031),
    @(3.805606016746037,1.494310085590225,2.446077621900006),
    @(3.119575682053099,4.113574073864802,-0.2212830832371099),
    @(-2.658262938719541,2.234168899976305,2.193659129509576),
    @(-2.523965861247135,9.674394070185159,1.710836189526775),
    @(-0.4347282877335233,7.832854989858879,-0.6348883940623381),
    @(2.215493936378243,-6.333109965691214,-3.700846394667294),
    @(3.503521535946755,-4.84552946755515,-5.782033076653083),
    @(0.5721309762734461,-0.04020348672683609,-0.7832388446880931),
    @(-4.998460626602359,-7.086485686898527,0.6293182730674433)
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
